$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MyApplication")

# Push the existing "VM / artifact0304 / ... / SouthZone" row (row 3) down to row 4,
# carrying its values and formatting along with it.
$ws.Range("A3:L3").Copy($ws.Range("A4:L4"))

# Seed the new row 3 from row 2's formatting/content (both are "container" entries),
# so the new row picks up the same styles without minting new style records.
$ws.Range("A2:L2").Copy($ws.Range("A3:L3"))

# The ApplicationNameFrDeProvision / ApplicationNmFrDeboard cells (K3/L3) on this new
# row use the narrower style from column I rather than row 2's K/L style, so re-stamp
# just the formatting for those two cells.
$ws.Range("I2").Copy()
$ws.Range("K3:L3").PasteSpecial(-4122)

# Now fill in the new application/artifact identifiers for this validation row.
$ws.Range("B3").Value = "artifact0311"
$ws.Range("E3").Value = "appdemo213"
$ws.Range("I3").Value = "appdemo213"
$ws.Range("K3").Value = "appdemo213"
$ws.Range("L3").Value = "appdemo213"

# Match the saved selection/scroll position recorded for this sheet.
$null = $ws.Range("K13").Select()
